$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.311.57"
$ws.Range("E2").Value = "'  +2.57%  "
$ws.Range("D3").Value = "'2.655.41"
$ws.Range("E3").Value = "'  +1.87%  "
$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("D5").Value = "'601.46"
$ws.Range("E5").Value = "'  +1.10%  "
$ws.Range("D6").Value = "'157.29"
$ws.Range("E6").Value = "'  +4.12%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "'  +0.03%  "
$ws.Range("D8").Value = "'0.593"
$ws.Range("E8").Value = "'  +0.95%  "
$ws.Range("E9").Value = "'  +6.90%  "
$ws.Range("D10").Value = "'0.402"
$ws.Range("E10").Value = "'  +4.40%  "
$ws.Range("E11").Value = "'  +2.66%  "
$ws.Range("E12").Value = "'  +1.47%  "
$ws.Range("D13").Value = "'29.22"
$ws.Range("E13").Value = "'  +5.02%  "
$ws.Range("D14").Value = "'3.134.47"
$ws.Range("E14").Value = "'  +1.80%  "
$ws.Range("B15").Value = "'ShibaInu"
$ws.Range("C15").Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "'0.0000174"
$ws.Range("E15").Value = "'  +13.19%  "
$ws.Range("B16").Value = "'WrappedBTC"
$ws.Range("C16").Value = "'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "'65.220.59"
$ws.Range("E16").Value = "'  +2.67%  "
$ws.Range("D17").Value = "'2.640.15"
$ws.Range("E17").Value = "'  +0.41%  "
$ws.Range("D18").Value = "'12.61"
$ws.Range("E18").Value = "'  +1.74%  "
$ws.Range("E19").Value = "'  +2.47%  "
$ws.Range("D20").Value = "'354.58"
$ws.Range("E20").Value = "'  +1.77%  "
$ws.Range("D21").Value = "'7.30"
$ws.Range("E21").Value = "'  +5.99%  "
$ws.Range("E22").Value = "'  +0.24%  "
$ws.Range("D23").Value = "'68.24"
$ws.Range("E23").Value = "'  +0.92%  "
$ws.Range("D24").Value = "'1.71"
$ws.Range("E24").Value = "'  -0.12%  "
$ws.Range("E25").Value = "'  +2.81%  "
$ws.Range("E26").Value = "'  -1.43%  "
$ws.Range("D27").Value = "'8.35"
$ws.Range("E27").Value = "'  +4.02%  "
$ws.Range("D28").Value = "'0.165"
$ws.Range("E28").Value = "'  +2.36%  "
$ws.Range("B29").Value = "'Binance-PegBSC-USD"
$ws.Range("C29").Value = "'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "'  -0.29%  "
$ws.Range("B30").Value = "'Bittensor"
$ws.Range("C30").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D30").Value = "'538.31"
$ws.Range("E30").Value = "'  -2.99%  "
$ws.Range("E31").Value = "'  +8.91%  "
$ws.Range("D32").Value = "'2.08"
$ws.Range("E32").Value = "'  +1.05%  "
$ws.Range("D33").Value = "'1.83"
$ws.Range("E33").Value = "'  +4.53%  "
$ws.Range("E34").Value = "'  +11.42%  "
$ws.Range("D35").Value = "'6.48"
$ws.Range("E35").Value = "'  +4.85%  "
$ws.Range("D36").Value = "'0.428"
$ws.Range("E36").Value = "'  +2.96%  "
$ws.Range("D37").Value = "'2.05"
$ws.Range("E37").Value = "'  +5.94%  "
$ws.Range("D38").Value = "'165.36"
$ws.Range("E38").Value = "'  -0.91%  "
$ws.Range("D39").Value = "'20.24"
$ws.Range("E39").Value = "'  +3.29%  "
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "'  +0.10%  "
$ws.Range("E41").Value = "'  -0.12%  "
$ws.Range("D42").Value = "'168.96"
$ws.Range("E42").Value = "'  +1.34%  "
$ws.Range("D43").Value = "'41.70"
$ws.Range("E43").Value = "'  +4.99%  "
$ws.Range("E44").Value = "'  +4.00%  "
$ws.Range("D45").Value = "'0.0611"
$ws.Range("E45").Value = "'  +3.92%  "
$ws.Range("D46").Value = "'23.44"
$ws.Range("E46").Value = "'  +6.58%  "
$ws.Range("D47").Value = "'2.25"
$ws.Range("E47").Value = "'  +11.15%  "
$ws.Range("D48").Value = "'0.650"
$ws.Range("E48").Value = "'  +2.78%  "
$ws.Range("D49").Value = "'0.0252"
$ws.Range("E49").Value = "'  -0.18%  "
$ws.Range("D50").Value = "'0.0984"
$ws.Range("E50").Value = "'  +1.87%  "
$ws.Range("D51").Value = "'19.50"
$ws.Range("E51").Value = "'  +1.44%  "

Write-Host "Updated cryptos worksheet"
